$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 10003
$ws.Range("J7").Value = 10003
$ws.Range("L7").Value = 10003
$ws.Range("N7").Value = -10227

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H14").Value = 10003
$ws.Range("J14").Value = 10003
$ws.Range("L14").Value = 10003
$ws.Range("N14").Value = -10385

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 12749.5
$ws.Range("I21").Value = 10713.714
$ws.Range("K21").Value = 10713.714
$ws.Range("M21").Value = -10245.714

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 12749.5
$ws.Range("I23").Value = 10713.714
$ws.Range("K23").Value = 10713.714
$ws.Range("M23").Value = -10479.714

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 18241.334
$ws.Range("I34").Value = 18241.334
$ws.Range("K34").Value = 18241.334
$ws.Range("M34").Value = -18038.334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H36").Value = 18241.334
$ws.Range("I36").Value = 18241.334
$ws.Range("K36").Value = 18241.334
$ws.Range("M36").Value = -17526.334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1041.9412
$ws.Range("I38").Value = 92.75
$ws.Range("J38").Value = 3320
$ws.Range("K38").Value = 278.25
$ws.Range("L38").Value = 9960
$ws.Range("M38").Value = 93.75
$ws.Range("N38").Value = -10704

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1243.6875
$ws.Range("I58").Value = 133.22223
$ws.Range("J58").Value = 2671.4285
$ws.Range("K58").Value = 399.66669
$ws.Range("L58").Value = 8014.2855
$ws.Range("M58").Value = -249.66669
$ws.Range("N58").Value = -8314.2855

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 847.2439000000001
$ws.Range("I135").Value = 573.15625
$ws.Range("K135").Value = 5158.40625
$ws.Range("M135").Value = -2623.40625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1923.2
$ws.Range("I137").Value = 1815.1578
$ws.Range("J137").Value = 2109.818
$ws.Range("K137").Value = 5445.4734
$ws.Range("L137").Value = 6329.454000000001
$ws.Range("M137").Value = -2895.4734
$ws.Range("N137").Value = -11429.454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2268.2856
$ws.Range("I61").Value = 2316.3333
$ws.Range("J61").Value = 1980
$ws.Range("K61").Value = 2316.3333
$ws.Range("L61").Value = 1980
$ws.Range("M61").Value = -2104.3333
$ws.Range("N61").Value = -2404

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2268.2856
$ws.Range("I136").Value = 2316.3333
$ws.Range("J136").Value = 1980
$ws.Range("K136").Value = 6948.999899999999
$ws.Range("L136").Value = 5940
$ws.Range("M136").Value = -4398.999899999999
$ws.Range("N136").Value = -11040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25003044
$ws.Range("I31").Value = 38463856
$ws.Range("K31").Value = 38463856
$ws.Range("M31").Value = -38463561

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 25003044
$ws.Range("I34").Value = 38463856
$ws.Range("K34").Value = 38463856
$ws.Range("M34").Value = -38463654

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1344.5143
$ws.Range("I58").Value = 1295.5625
$ws.Range("J58").Value = 1866.6666
$ws.Range("K58").Value = 1295.5625
$ws.Range("L58").Value = 1866.6666
$ws.Range("M58").Value = -1092.5625
$ws.Range("N58").Value = -2272.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1344.5143
$ws.Range("I136").Value = 1295.5625
$ws.Range("J136").Value = 1866.6666
$ws.Range("K136").Value = 3886.6875
$ws.Range("L136").Value = 5599.9998
$ws.Range("M136").Value = -1336.6875
$ws.Range("N136").Value = -10699.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 2338.2
$ws.Range("I25").Value = 395.5
$ws.Range("K25").Value = 1186.5
$ws.Range("M25").Value = -1017.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H30").Value = 2338.2
$ws.Range("I30").Value = 395.5
$ws.Range("K30").Value = 1186.5
$ws.Range("M30").Value = -1084.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1050
$ws.Range("I34").Value = 200
$ws.Range("J34").Value = 1333.3334
$ws.Range("K34").Value = 600
$ws.Range("L34").Value = 4000.0002
$ws.Range("M34").Value = -516
$ws.Range("N34").Value = -4168.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 5683.6665
$ws.Range("I55").Value = 700
$ws.Range("K55").Value = 2100
$ws.Range("M55").Value = -1923

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 2010.3334
$ws.Range("J98").Value = 1281.6666
$ws.Range("L98").Value = 3844.9998
$ws.Range("N98").Value = -6840.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 6455
$ws.Range("I124").Value = 932.5
$ws.Range("K124").Value = 2797.5
$ws.Range("M124").Value = 2112.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1516362
$ws.Range("I129").Value = 452.41666
$ws.Range("J129").Value = 2382596
$ws.Range("K129").Value = 1357.24998
$ws.Range("L129").Value = 7147788
$ws.Range("M129").Value = 3642.75002
$ws.Range("N129").Value = -7157788

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 2182
$ws.Range("I130").Value = 806.6667
$ws.Range("J130").Value = 2771.4285
$ws.Range("K130").Value = 2420.0001
$ws.Range("L130").Value = 8314.2855
$ws.Range("M130").Value = 2599.9999
$ws.Range("N130").Value = -18354.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 864.9299999999999
$ws.Range("J131").Value = 885.29474
$ws.Range("L131").Value = 2655.88422
$ws.Range("N131").Value = -12735.88422

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 20000
$ws.Range("J39").Value = 20000
$ws.Range("L39").Value = 20000
$ws.Range("N39").Value = -21064

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 11294.417
$ws.Range("J123").Value = 11294.417
$ws.Range("L123").Value = 11294.417
$ws.Range("N123").Value = -16194.417

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 518.6667
$ws.Range("I22").Value = 450
$ws.Range("J22").Value = 656
$ws.Range("K22").Value = 450
$ws.Range("L22").Value = 656
$ws.Range("M22").Value = -155
$ws.Range("N22").Value = -1246

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 518.6667
$ws.Range("I27").Value = 450
$ws.Range("J27").Value = 656
$ws.Range("K27").Value = 450
$ws.Range("L27").Value = 656
$ws.Range("M27").Value = -343
$ws.Range("N27").Value = -870

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 17386
$ws.Range("I61").Value = 23150.5
$ws.Range("K61").Value = 23150.5
$ws.Range("M61").Value = -22948.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 17386
$ws.Range("I113").Value = 23150.5
$ws.Range("K113").Value = 23150.5
$ws.Range("M113").Value = -20980.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5168.7144
$ws.Range("I132").Value = 4903.1562
$ws.Range("J132").Value = 8001.3335
$ws.Range("K132").Value = 14709.4686
$ws.Range("L132").Value = 24004.0005
$ws.Range("M132").Value = -12179.4686
$ws.Range("N132").Value = -29064.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3284.1924
$ws.Range("I136").Value = 3499.9556
$ws.Range("K136").Value = 10499.8668
$ws.Range("M136").Value = -7949.8668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4429
$ws.Range("I132").Value = 4329.5713
$ws.Range("J132").Value = 4528.4287
$ws.Range("K132").Value = 12988.7139
$ws.Range("L132").Value = 13585.2861
$ws.Range("M132").Value = -10458.7139
$ws.Range("N132").Value = -18645.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1448.081
$ws.Range("I136").Value = 1591.1034
$ws.Range("J136").Value = 929.625
$ws.Range("K136").Value = 4773.3102
$ws.Range("L136").Value = 2788.875
$ws.Range("M136").Value = -2223.3102
$ws.Range("N136").Value = -7888.875
